$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they remain stored as text
# (matching the source data, which is all text).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"


# Row 2
$ws.Range("D2").Value = '67.752.85'
$ws.Range("E2").Value = '  +6.34%  '

# Row 3
$ws.Range("D3").Value = '3.711.83'
$ws.Range("E3").Value = '  +6.55%  '

# Row 4
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.65%  '

# Row 5
$ws.Range("D5").Value = '426.98'
$ws.Range("E5").Value = '  +2.89%  '

# Row 6
$ws.Range("D6").Value = '130.90'
$ws.Range("E6").Value = '  +0.85%  '

# Row 7
$ws.Range("D7").Value = '3.703.51'
$ws.Range("E7").Value = '  +6.46%  '

# Row 8
$ws.Range("D8").Value = '0.645'
$ws.Range("E8").Value = '  +2.16%  '

# Row 9
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.04%  '

# Row 10
$ws.Range("D10").Value = '0.768'
$ws.Range("E10").Value = '  +2.74%  '

# Row 11
$ws.Range("D11").Value = '0.182'
$ws.Range("E11").Value = '  +14.89%  '

# Row 12
$ws.Range("D12").Value = '0.0000388'
$ws.Range("E12").Value = '  +70.39%  '

# Row 13
$ws.Range("D13").Value = '42.76'
$ws.Range("E13").Value = '  +0.76%  '

# Row 14
$ws.Range("D14").Value = '10.00'
$ws.Range("E14").Value = '  +2.55%  '

# Row 15
$ws.Range("D15").Value = '4.300.16'
$ws.Range("E15").Value = '  +6.53%  '

# Row 16
$ws.Range("E16").Value = '  +0.14%  '

# Row 17
$ws.Range("D17").Value = '20.63'
$ws.Range("E17").Value = '  +1.64%  '

# Row 18
$ws.Range("D18").Value = '3.720.72'
$ws.Range("E18").Value = '  +6.97%  '

# Row 19
$ws.Range("E19").Value = '  +3.20%  '

# Row 20
$ws.Range("D20").Value = '12.90'
$ws.Range("E20").Value = '  +3.57%  '

# Row 21
$ws.Range("D21").Value = '67.826.44'
$ws.Range("E21").Value = '  +6.57%  '

# Row 22
$ws.Range("D22").Value = '450.20'
$ws.Range("E22").Value = '  -2.35%  '

# Row 23
$ws.Range("E23").Value = '  +15.33%  '

# Row 24
$ws.Range("D24").Value = '89.65'
$ws.Range("E24").Value = '  -0.73%  '

# Row 25
$ws.Range("D25").Value = '3.14'
$ws.Range("E25").Value = '  -4.33%  '

# Row 26
$ws.Range("D26").Value = '38.58'
$ws.Range("E26").Value = '  +15.02%  '

# Row 27
$ws.Range("D27").Value = '10.35'
$ws.Range("E27").Value = '  +1.57%  '

# Row 28
$ws.Range("D28").Value = '3.32'
$ws.Range("E28").Value = '  +0.04%  '

# Row 29
$ws.Range("D29").Value = '4.97'
$ws.Range("E29").Value = '  +4.31%  '

# Row 30
$ws.Range("D30").Value = '2.80'
$ws.Range("E30").Value = '  +4.37%  '

# Row 31
$ws.Range("D31").Value = '12.51'
$ws.Range("E31").Value = '  +0.71%  '

# Row 32
$ws.Range("E32").Value = '  +8.04%  '

# Row 33
$ws.Range("E33").Value = '  -4.00%  '

# Row 34
$ws.Range("E34").Value = '  -3.75%  '

# Row 35
$ws.Range("D35").Value = '40.43'
$ws.Range("E35").Value = '  +0.72%  '

# Row 36
$ws.Range("E36").Value = '  -0.14%  '

# Row 37
$ws.Range("D37").Value = '56.51'
$ws.Range("E37").Value = '  -3.48%  '

# Row 38
$ws.Range("D38").Value = '0.0491'
$ws.Range("E38").Value = '  +1.00%  '

# Row 39
$ws.Range("D39").Value = '0.0₃0716'
$ws.Range("E39").Value = '  +11.17%  '

# Row 40
$ws.Range("D40").Value = '3.02'
$ws.Range("E40").Value = '  +28.69%  '

# Row 41
$ws.Range("D41").Value = '0.148'
$ws.Range("E41").Value = '  +7.33%  '

# Row 42
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.06%  '

# Row 43
$ws.Range("D43").Value = '3.44'
$ws.Range("E43").Value = '  +2.98%  '

# Row 44
$ws.Range("D44").Value = '147.23'
$ws.Range("E44").Value = '  +0.81%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '2.92'
$ws.Range("E45").Value = '  -5.86%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '2.68'
$ws.Range("E46").Value = '  -4.99%  '

# Row 47
$ws.Range("D47").Value = '26.62'
$ws.Range("E47").Value = '  +22.44%  '

# Row 48
$ws.Range("D48").Value = '2.09'
$ws.Range("E48").Value = '  +4.36%  '

# Row 49
$ws.Range("D49").Value = '4.33'
$ws.Range("E49").Value = '  -4.64%  '

# Row 50
$ws.Range("D50").Value = '0.306'
$ws.Range("E50").Value = '  -3.34%  '

# Row 51
$ws.Range("E51").Value = '  +15.45%  '
